# The scraper now also records height/weight for each player. These two new
# columns are inserted before the existing "fantasy points" column, so the
# old column E ("fantasy points") data moves over to column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing "fantasy points" values (column E, rows 2-17) before
# we overwrite anything, so we can re-write them untouched into column G.
$fantasyPoints = @{}
for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints[$r] = $ws.Cells.Item($r, 5).Value()
}

# New header labels for the inserted columns, matching the style used by the
# other header cells (bold, bordered, centered) by copying D1's style.
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Fill in the height/weight values, and re-write the fantasy points values
# into their new home (column G) unchanged.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.25
    $ws.Cells.Item($r, 6).Value = 255
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r]
}
